$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly re-shuffle of Fecha/Volumen/Precio fields across existing rows
# (D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg)

$ws.Cells.Item(2, 4).Value = 44664   # D2 Fecha
$ws.Cells.Item(2, 10).Value = 50  # J2 Volumen
$ws.Cells.Item(2, 11).Value = 11000  # K2 Precio minimo
$ws.Cells.Item(2, 12).Value = 12000  # L2 Precio maximo
$ws.Cells.Item(2, 13).Value = 11600  # M2 Precio promedio ponderado
$ws.Cells.Item(2, 16).Value = 892  # P2 Precio $/Kg

$ws.Cells.Item(3, 4).Value = 44838   # D3 Fecha
$ws.Cells.Item(3, 10).Value = 40  # J3 Volumen
$ws.Cells.Item(3, 11).Value = 14000  # K3 Precio minimo
$ws.Cells.Item(3, 12).Value = 15000  # L3 Precio maximo
$ws.Cells.Item(3, 13).Value = 14500  # M3 Precio promedio ponderado
$ws.Cells.Item(3, 16).Value = 1115  # P3 Precio $/Kg

$ws.Cells.Item(4, 4).Value = 44810   # D4 Fecha
$ws.Cells.Item(4, 10).Value = 50  # J4 Volumen
$ws.Cells.Item(4, 11).Value = 11000  # K4 Precio minimo
$ws.Cells.Item(4, 12).Value = 12000  # L4 Precio maximo
$ws.Cells.Item(4, 13).Value = 11600  # M4 Precio promedio ponderado
$ws.Cells.Item(4, 16).Value = 892  # P4 Precio $/Kg

$ws.Cells.Item(5, 4).Value = 44377   # D5 Fecha
$ws.Cells.Item(5, 10).Value = 40  # J5 Volumen
$ws.Cells.Item(5, 11).Value = 14000  # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 15000  # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 14500  # M5 Precio promedio ponderado
$ws.Cells.Item(5, 16).Value = 1115  # P5 Precio $/Kg

$ws.Cells.Item(6, 4).Value = 44320   # D6 Fecha
$ws.Cells.Item(6, 10).Value = 50  # J6 Volumen
$ws.Cells.Item(6, 11).Value = 26000  # K6 Precio minimo
$ws.Cells.Item(6, 12).Value = 28000  # L6 Precio maximo
$ws.Cells.Item(6, 13).Value = 26800  # M6 Precio promedio ponderado
$ws.Cells.Item(6, 16).Value = 2062  # P6 Precio $/Kg

$ws.Cells.Item(7, 4).Value = 44523   # D7 Fecha
$ws.Cells.Item(7, 10).Value = 40  # J7 Volumen
$ws.Cells.Item(7, 11).Value = 15000  # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 16000  # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 15500  # M7 Precio promedio ponderado
$ws.Cells.Item(7, 16).Value = 1192  # P7 Precio $/Kg

$ws.Cells.Item(8, 4).Value = 44719   # D8 Fecha
$ws.Cells.Item(8, 10).Value = 50  # J8 Volumen
$ws.Cells.Item(8, 11).Value = 13000  # K8 Precio minimo
$ws.Cells.Item(8, 12).Value = 14000  # L8 Precio maximo
$ws.Cells.Item(8, 13).Value = 13400  # M8 Precio promedio ponderado
$ws.Cells.Item(8, 16).Value = 1031  # P8 Precio $/Kg

$ws.Cells.Item(9, 4).Value = 44462   # D9 Fecha
$ws.Cells.Item(9, 10).Value = 60  # J9 Volumen
$ws.Cells.Item(9, 11).Value = 14000  # K9 Precio minimo
$ws.Cells.Item(9, 12).Value = 15000  # L9 Precio maximo
$ws.Cells.Item(9, 13).Value = 14500  # M9 Precio promedio ponderado
$ws.Cells.Item(9, 16).Value = 1115  # P9 Precio $/Kg

$ws.Cells.Item(10, 4).Value = 44755   # D10 Fecha
$ws.Cells.Item(10, 10).Value = 40  # J10 Volumen
$ws.Cells.Item(10, 11).Value = 14000  # K10 Precio minimo
$ws.Cells.Item(10, 12).Value = 15000  # L10 Precio maximo
$ws.Cells.Item(10, 13).Value = 14500  # M10 Precio promedio ponderado
$ws.Cells.Item(10, 16).Value = 1115  # P10 Precio $/Kg

$ws.Cells.Item(11, 4).Value = 44467   # D11 Fecha
$ws.Cells.Item(11, 10).Value = 100  # J11 Volumen
$ws.Cells.Item(11, 11).Value = 13000  # K11 Precio minimo
$ws.Cells.Item(11, 12).Value = 14000  # L11 Precio maximo
$ws.Cells.Item(11, 13).Value = 13500  # M11 Precio promedio ponderado
$ws.Cells.Item(11, 16).Value = 1038  # P11 Precio $/Kg

$ws.Cells.Item(12, 4).Value = 44510   # D12 Fecha
$ws.Cells.Item(12, 10).Value = 40  # J12 Volumen
$ws.Cells.Item(12, 11).Value = 15000  # K12 Precio minimo
$ws.Cells.Item(12, 12).Value = 16000  # L12 Precio maximo
$ws.Cells.Item(12, 13).Value = 15500  # M12 Precio promedio ponderado
$ws.Cells.Item(12, 16).Value = 1192  # P12 Precio $/Kg

$ws.Cells.Item(13, 4).Value = 44691   # D13 Fecha
$ws.Cells.Item(13, 10).Value = 100  # J13 Volumen
$ws.Cells.Item(13, 11).Value = 12000  # K13 Precio minimo
$ws.Cells.Item(13, 12).Value = 13000  # L13 Precio maximo
$ws.Cells.Item(13, 13).Value = 12500  # M13 Precio promedio ponderado
$ws.Cells.Item(13, 16).Value = 962  # P13 Precio $/Kg

$ws.Cells.Item(14, 4).Value = 44383   # D14 Fecha
$ws.Cells.Item(14, 10).Value = 50  # J14 Volumen
$ws.Cells.Item(14, 11).Value = 15000  # K14 Precio minimo
$ws.Cells.Item(14, 12).Value = 16000  # L14 Precio maximo
$ws.Cells.Item(14, 13).Value = 15400  # M14 Precio promedio ponderado
$ws.Cells.Item(14, 16).Value = 1185  # P14 Precio $/Kg

$ws.Cells.Item(15, 4).Value = 44433   # D15 Fecha
$ws.Cells.Item(15, 10).Value = 100  # J15 Volumen
$ws.Cells.Item(15, 11).Value = 13000  # K15 Precio minimo
$ws.Cells.Item(15, 12).Value = 14000  # L15 Precio maximo
$ws.Cells.Item(15, 13).Value = 13500  # M15 Precio promedio ponderado
$ws.Cells.Item(15, 16).Value = 1038  # P15 Precio $/Kg

$ws.Cells.Item(16, 4).Value = 44435   # D16 Fecha
$ws.Cells.Item(16, 10).Value = 100  # J16 Volumen
$ws.Cells.Item(16, 11).Value = 13000  # K16 Precio minimo
$ws.Cells.Item(16, 12).Value = 14000  # L16 Precio maximo
$ws.Cells.Item(16, 13).Value = 13500  # M16 Precio promedio ponderado
$ws.Cells.Item(16, 16).Value = 1038  # P16 Precio $/Kg

$ws.Cells.Item(17, 4).Value = 44313   # D17 Fecha
$ws.Cells.Item(17, 10).Value = 50  # J17 Volumen
$ws.Cells.Item(17, 11).Value = 25000  # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 26000  # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 25600  # M17 Precio promedio ponderado
$ws.Cells.Item(17, 16).Value = 1969  # P17 Precio $/Kg

$ws.Cells.Item(18, 4).Value = 44159   # D18 Fecha
$ws.Cells.Item(18, 10).Value = 60  # J18 Volumen
$ws.Cells.Item(18, 11).Value = 30000  # K18 Precio minimo
$ws.Cells.Item(18, 12).Value = 32000  # L18 Precio maximo
$ws.Cells.Item(18, 13).Value = 31000  # M18 Precio promedio ponderado
$ws.Cells.Item(18, 16).Value = 2385  # P18 Precio $/Kg

$ws.Cells.Item(19, 4).Value = 44316   # D19 Fecha
$ws.Cells.Item(19, 10).Value = 50  # J19 Volumen
$ws.Cells.Item(19, 11).Value = 27000  # K19 Precio minimo
$ws.Cells.Item(19, 12).Value = 28000  # L19 Precio maximo
$ws.Cells.Item(19, 13).Value = 27400  # M19 Precio promedio ponderado
$ws.Cells.Item(19, 16).Value = 2108  # P19 Precio $/Kg

$ws.Cells.Item(21, 4).Value = 44308   # D21 Fecha
$ws.Cells.Item(21, 10).Value = 50  # J21 Volumen
$ws.Cells.Item(21, 11).Value = 26000  # K21 Precio minimo
$ws.Cells.Item(21, 12).Value = 27000  # L21 Precio maximo
$ws.Cells.Item(21, 13).Value = 26400  # M21 Precio promedio ponderado
$ws.Cells.Item(21, 16).Value = 2031  # P21 Precio $/Kg

$ws.Cells.Item(22, 4).Value = 44355   # D22 Fecha
$ws.Cells.Item(22, 10).Value = 60  # J22 Volumen
$ws.Cells.Item(22, 11).Value = 18000  # K22 Precio minimo
$ws.Cells.Item(22, 12).Value = 20000  # L22 Precio maximo
$ws.Cells.Item(22, 13).Value = 19000  # M22 Precio promedio ponderado
$ws.Cells.Item(22, 16).Value = 1462  # P22 Precio $/Kg

$ws.Cells.Item(23, 4).Value = 44708   # D23 Fecha
$ws.Cells.Item(23, 10).Value = 50  # J23 Volumen
$ws.Cells.Item(23, 11).Value = 13000  # K23 Precio minimo
$ws.Cells.Item(23, 12).Value = 14000  # L23 Precio maximo
$ws.Cells.Item(23, 13).Value = 13600  # M23 Precio promedio ponderado
$ws.Cells.Item(23, 16).Value = 1046  # P23 Precio $/Kg

$ws.Cells.Item(24, 4).Value = 44488   # D24 Fecha
$ws.Cells.Item(24, 10).Value = 40  # J24 Volumen
$ws.Cells.Item(24, 11).Value = 16000  # K24 Precio minimo
$ws.Cells.Item(24, 12).Value = 17000  # L24 Precio maximo
$ws.Cells.Item(24, 13).Value = 16500  # M24 Precio promedio ponderado
$ws.Cells.Item(24, 16).Value = 1269  # P24 Precio $/Kg

$ws.Cells.Item(25, 4).Value = 44610   # D25 Fecha
$ws.Cells.Item(25, 10).Value = 50  # J25 Volumen
$ws.Cells.Item(25, 11).Value = 17000  # K25 Precio minimo
$ws.Cells.Item(25, 12).Value = 18000  # L25 Precio maximo
$ws.Cells.Item(25, 13).Value = 17400  # M25 Precio promedio ponderado
$ws.Cells.Item(25, 16).Value = 1338  # P25 Precio $/Kg

$ws.Cells.Item(26, 4).Value = 44327   # D26 Fecha
$ws.Cells.Item(26, 10).Value = 50  # J26 Volumen
$ws.Cells.Item(26, 11).Value = 24000  # K26 Precio minimo
$ws.Cells.Item(26, 12).Value = 25000  # L26 Precio maximo
$ws.Cells.Item(26, 13).Value = 24400  # M26 Precio promedio ponderado
$ws.Cells.Item(26, 16).Value = 1877  # P26 Precio $/Kg

$ws.Cells.Item(27, 4).Value = 44362   # D27 Fecha
$ws.Cells.Item(27, 10).Value = 40  # J27 Volumen
$ws.Cells.Item(27, 11).Value = 15000  # K27 Precio minimo
$ws.Cells.Item(27, 12).Value = 16000  # L27 Precio maximo
$ws.Cells.Item(27, 13).Value = 15500  # M27 Precio promedio ponderado
$ws.Cells.Item(27, 16).Value = 1192  # P27 Precio $/Kg

$ws.Cells.Item(28, 4).Value = 44503   # D28 Fecha
$ws.Cells.Item(28, 10).Value = 35  # J28 Volumen
$ws.Cells.Item(28, 11).Value = 15000  # K28 Precio minimo
$ws.Cells.Item(28, 12).Value = 16000  # L28 Precio maximo
$ws.Cells.Item(28, 13).Value = 15429  # M28 Precio promedio ponderado
$ws.Cells.Item(28, 16).Value = 1187  # P28 Precio $/Kg

$ws.Cells.Item(29, 4).Value = 44509   # D29 Fecha
$ws.Cells.Item(29, 10).Value = 100  # J29 Volumen
$ws.Cells.Item(29, 11).Value = 15000  # K29 Precio minimo
$ws.Cells.Item(29, 12).Value = 16000  # L29 Precio maximo
$ws.Cells.Item(29, 13).Value = 15500  # M29 Precio promedio ponderado
$ws.Cells.Item(29, 16).Value = 1192  # P29 Precio $/Kg

$ws.Cells.Item(30, 4).Value = 44775   # D30 Fecha
$ws.Cells.Item(30, 10).Value = 20  # J30 Volumen
$ws.Cells.Item(30, 11).Value = 12000  # K30 Precio minimo
$ws.Cells.Item(30, 12).Value = 13000  # L30 Precio maximo
$ws.Cells.Item(30, 13).Value = 12500  # M30 Precio promedio ponderado
$ws.Cells.Item(30, 16).Value = 962  # P30 Precio $/Kg

$ws.Cells.Item(31, 4).Value = 44782   # D31 Fecha
$ws.Cells.Item(31, 10).Value = 40  # J31 Volumen
$ws.Cells.Item(31, 11).Value = 13000  # K31 Precio minimo
$ws.Cells.Item(31, 12).Value = 14000  # L31 Precio maximo
$ws.Cells.Item(31, 13).Value = 13500  # M31 Precio promedio ponderado
$ws.Cells.Item(31, 16).Value = 1038  # P31 Precio $/Kg

$ws.Cells.Item(33, 4).Value = 44705   # D33 Fecha
$ws.Cells.Item(33, 10).Value = 50  # J33 Volumen
$ws.Cells.Item(33, 11).Value = 10000  # K33 Precio minimo
$ws.Cells.Item(33, 12).Value = 11000  # L33 Precio maximo
$ws.Cells.Item(33, 13).Value = 10400  # M33 Precio promedio ponderado
$ws.Cells.Item(33, 16).Value = 800  # P33 Precio $/Kg

$ws.Cells.Item(34, 4).Value = 44777   # D34 Fecha
$ws.Cells.Item(34, 10).Value = 25  # J34 Volumen
$ws.Cells.Item(34, 11).Value = 13000  # K34 Precio minimo
$ws.Cells.Item(34, 12).Value = 14000  # L34 Precio maximo
$ws.Cells.Item(34, 13).Value = 13600  # M34 Precio promedio ponderado
$ws.Cells.Item(34, 16).Value = 1046  # P34 Precio $/Kg

$ws.Cells.Item(35, 4).Value = 44474   # D35 Fecha
$ws.Cells.Item(35, 10).Value = 40  # J35 Volumen
$ws.Cells.Item(35, 11).Value = 13000  # K35 Precio minimo
$ws.Cells.Item(35, 12).Value = 14000  # L35 Precio maximo
$ws.Cells.Item(35, 13).Value = 13500  # M35 Precio promedio ponderado
$ws.Cells.Item(35, 16).Value = 1038  # P35 Precio $/Kg

$ws.Cells.Item(36, 4).Value = 44761   # D36 Fecha
$ws.Cells.Item(36, 10).Value = 25  # J36 Volumen
$ws.Cells.Item(36, 11).Value = 14000  # K36 Precio minimo
$ws.Cells.Item(36, 12).Value = 15000  # L36 Precio maximo
$ws.Cells.Item(36, 13).Value = 14400  # M36 Precio promedio ponderado
$ws.Cells.Item(36, 16).Value = 1108  # P36 Precio $/Kg

$ws.Cells.Item(37, 4).Value = 44425   # D37 Fecha
$ws.Cells.Item(37, 10).Value = 60  # J37 Volumen
$ws.Cells.Item(37, 11).Value = 14000  # K37 Precio minimo
$ws.Cells.Item(37, 12).Value = 15000  # L37 Precio maximo
$ws.Cells.Item(37, 13).Value = 14500  # M37 Precio promedio ponderado
$ws.Cells.Item(37, 16).Value = 1115  # P37 Precio $/Kg

$ws.Cells.Item(38, 4).Value = 44264   # D38 Fecha
$ws.Cells.Item(38, 10).Value = 40  # J38 Volumen
$ws.Cells.Item(38, 11).Value = 30000  # K38 Precio minimo
$ws.Cells.Item(38, 12).Value = 32000  # L38 Precio maximo
$ws.Cells.Item(38, 13).Value = 31000  # M38 Precio promedio ponderado
$ws.Cells.Item(38, 16).Value = 2385  # P38 Precio $/Kg

$ws.Cells.Item(39, 4).Value = 44453   # D39 Fecha
$ws.Cells.Item(39, 10).Value = 50  # J39 Volumen
$ws.Cells.Item(39, 11).Value = 14000  # K39 Precio minimo
$ws.Cells.Item(39, 12).Value = 15000  # L39 Precio maximo
$ws.Cells.Item(39, 13).Value = 14600  # M39 Precio promedio ponderado
$ws.Cells.Item(39, 16).Value = 1123  # P39 Precio $/Kg

$ws.Cells.Item(40, 4).Value = 44883   # D40 Fecha
$ws.Cells.Item(40, 10).Value = 60  # J40 Volumen
$ws.Cells.Item(40, 11).Value = 14000  # K40 Precio minimo
$ws.Cells.Item(40, 12).Value = 15000  # L40 Precio maximo
$ws.Cells.Item(40, 13).Value = 14500  # M40 Precio promedio ponderado
$ws.Cells.Item(40, 16).Value = 1115  # P40 Precio $/Kg

$ws.Cells.Item(41, 4).Value = 44813   # D41 Fecha
$ws.Cells.Item(41, 10).Value = 50  # J41 Volumen
$ws.Cells.Item(41, 11).Value = 13000  # K41 Precio minimo
$ws.Cells.Item(41, 12).Value = 14000  # L41 Precio maximo
$ws.Cells.Item(41, 13).Value = 13400  # M41 Precio promedio ponderado
$ws.Cells.Item(41, 16).Value = 1031  # P41 Precio $/Kg

$ws.Cells.Item(42, 4).Value = 44819   # D42 Fecha
$ws.Cells.Item(42, 10).Value = 50  # J42 Volumen
$ws.Cells.Item(42, 11).Value = 13000  # K42 Precio minimo
$ws.Cells.Item(42, 12).Value = 14000  # L42 Precio maximo
$ws.Cells.Item(42, 13).Value = 13400  # M42 Precio promedio ponderado
$ws.Cells.Item(42, 16).Value = 1031  # P42 Precio $/Kg

$ws.Cells.Item(44, 4).Value = 44350   # D44 Fecha
$ws.Cells.Item(44, 10).Value = 40  # J44 Volumen
$ws.Cells.Item(44, 11).Value = 23000  # K44 Precio minimo
$ws.Cells.Item(44, 12).Value = 25000  # L44 Precio maximo
$ws.Cells.Item(44, 13).Value = 24000  # M44 Precio promedio ponderado
$ws.Cells.Item(44, 16).Value = 1846  # P44 Precio $/Kg
